# bug fix in 191
# Append new interview-history rows to AMSIN, BETA and AMS sheets, and
# correct a stray run-time value on AMS!B65.

$wb = $excel.ActiveWorkbook

# xlPasteFormats - used to clone a neighbouring cell's number format
# (in particular the "yyyy-mm-dd hh:mm:ss" style used by column B)
# onto newly written cells without Excel inventing a brand-new style.
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# AMSIN : new rows 98-104
# ---------------------------------------------------------------------
$wsAMSIN = $wb.Worksheets.Item("AMSIN")

$amsinData = @(
    @(98,  "2024-03-28", 45379.51146578704, "190masstrail",  155, 154, 1,  3.06),
    @(99,  "2024-03-29", 45380.40826674768, "190fstcycle",   155, 152, 3,  3.63),
    @(100, "2024-03-29", 45380.46256306713, "190scndcycle",  155, 155, 0,  2.97),
    @(101, "2024-04-01", 45383.38072798611, "190fnlrun",     155, 143, 12, 4.74),
    @(102, "2024-04-25", 45407.45873856481, "191trail",      155, 151, 4,  3.84),
    @(103, "2024-05-02", 45414.50755787037, "191fstcycle",   155, 149, 6,  3.54),
    @(104, "2024-05-03", 45415.35905292824, "191lstrun",     155, 155, 0,  3.16)
)

$amsinStyleRefA = $wsAMSIN.Cells.Item(97, 1)
$amsinStyleRefB = $wsAMSIN.Cells.Item(97, 2)

foreach ($row in $amsinData) {
    $r = $row[0]

    # Column A holds a date-look-alike string ("2024-03-28"). Mark the
    # cell as Text first so Excel doesn't silently convert it to a date
    # serial, then copy the regular style back onto it.
    $cA = $wsAMSIN.Cells.Item($r, 1)
    $cA.NumberFormat = "@"
    $cA.Value = $row[1]
    $cA.Style = $amsinStyleRefA.Style

    $cB = $wsAMSIN.Cells.Item($r, 2)
    $cB.Value = $row[2]
    $amsinStyleRefB.Copy()
    $cB.PasteSpecial($xlPasteFormats)

    $wsAMSIN.Cells.Item($r, 3).Value = $row[3]
    $wsAMSIN.Cells.Item($r, 4).Value = $row[4]
    $wsAMSIN.Cells.Item($r, 5).Value = $row[5]
    $wsAMSIN.Cells.Item($r, 6).Value = $row[6]
    $wsAMSIN.Cells.Item($r, 7).Value = $row[7]
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# BETA : new rows 45-46
# ---------------------------------------------------------------------
$wsBETA = $wb.Worksheets.Item("BETA")

$betaData = @(
    @(45, "2024-04-01", 45383.54779344908, "190betatest", 155, 153, 2, 4.05),
    @(46, "2024-05-03", 45415.62367593719, "191beta",     155, 155, 0, 3.35)
)

$betaStyleRefA = $wsBETA.Cells.Item(44, 1)
$betaStyleRefB = $wsBETA.Cells.Item(44, 2)

foreach ($row in $betaData) {
    $r = $row[0]

    $cA = $wsBETA.Cells.Item($r, 1)
    $cA.NumberFormat = "@"
    $cA.Value = $row[1]
    $cA.Style = $betaStyleRefA.Style

    $cB = $wsBETA.Cells.Item($r, 2)
    $cB.Value = $row[2]
    $betaStyleRefB.Copy()
    $cB.PasteSpecial($xlPasteFormats)

    $wsBETA.Cells.Item($r, 3).Value = $row[3]
    $wsBETA.Cells.Item($r, 4).Value = $row[4]
    $wsBETA.Cells.Item($r, 5).Value = $row[5]
    $wsBETA.Cells.Item($r, 6).Value = $row[6]
    $wsBETA.Cells.Item($r, 7).Value = $row[7]
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# AMS : fix row 65 (stray/odd B65 run time + missing explicit styling)
#       and append new row 66
# ---------------------------------------------------------------------
$wsAMS = $wb.Worksheets.Item("AMS")

# row 65 already holds the correct data, it only needs (a) the "normal"
# per-cell style applied (matching every other data row, e.g. row 64)
# and (b) the run-time value corrected.
$amsStyleRefA = $wsAMS.Cells.Item(64, 1)
$amsStyleRefC = $wsAMS.Cells.Item(64, 3)
$amsStyleRefD = $wsAMS.Cells.Item(64, 4)

$wsAMS.Cells.Item(65, 1).Style = $amsStyleRefA.Style
$wsAMS.Cells.Item(65, 3).Style = $amsStyleRefC.Style
$wsAMS.Range("D65:G65").Style = $amsStyleRefD.Style

$wsAMS.Cells.Item(65, 2).Value = 45359.75036765046

# new row 66
$amsData66 = @(66, "2024-04-01", 45383.86431619213, "190livee", 155, 155, 0, 3.53)

$amsStyleRefB = $wsAMS.Cells.Item(64, 2)

$r = $amsData66[0]
$cA = $wsAMS.Cells.Item($r, 1)
$cA.NumberFormat = "@"
$cA.Value = $amsData66[1]
$cA.Style = $amsStyleRefA.Style

$cB = $wsAMS.Cells.Item($r, 2)
$cB.Value = $amsData66[2]
$amsStyleRefB.Copy()
$cB.PasteSpecial($xlPasteFormats)

$wsAMS.Cells.Item($r, 3).Value = $amsData66[3]
$wsAMS.Cells.Item($r, 3).Style = $amsStyleRefC.Style
$wsAMS.Cells.Item($r, 4).Value = $amsData66[4]
$wsAMS.Cells.Item($r, 5).Value = $amsData66[5]
$wsAMS.Cells.Item($r, 6).Value = $amsData66[6]
$wsAMS.Cells.Item($r, 7).Value = $amsData66[7]
$wsAMS.Range("D66:G66").Style = $amsStyleRefD.Style

$excel.CutCopyMode = $false

Write-Output "Applied MASS_INTERVIEW_HISTORY_DATA updates (AMSIN +7 rows, BETA +2 rows, AMS row65 fix + 1 row)."
